$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (shifts old E -> F), bringing along formatting
# from the surrounding cells the same way Excel's UI "Insert Column" does.
$ws.Columns.Item(5).Insert()

# Resize the (now two) "35"-wide columns D and E.
# 34.285714285714285 (= 35 - 5/7) round-trips to a stored width of exactly 35.
$ws.Columns.Item(4).ColumnWidth = 34.285714285714285
$ws.Columns.Item(5).ColumnWidth = 34.285714285714285

# Fill in the new column's values.
$ws.Cells.Item(1,5).Value2 = "storeXpathCount"
$ws.Cells.Item(3,5).Value2 = '{"target":"xpath=//input","value":"input3"}'

# Adjust row heights: row 1 goes back to the default (no explicit height),
# rows 2 and 3 become 37.5.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).RowHeight = 37.5
$ws.Rows.Item(3).RowHeight = 37.5

# Update the active selection to match the saved view state.
$ws.Range("E8").Select() | Out-Null
